# Adds a "Почта" (email) column (K) to the sheet, populated with the
# parents' e-mail addresses that already live (as free text) inside
# column D ("Родители"), puts "-" where no e-mail is known, fills the
# one missing "Родители" cell (D13) with "-" as well, and turns four of
# the e-mail cells into real mailto: hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("K1").Value = "Почта"

# Plain-text e-mail / placeholder cells
$ws.Range("K2").Value  = " gtnhjdbx72@mail.ru"
$ws.Range("K3").Value  = "kalikina84@mail.ru"
$ws.Range("K4").Value  = "-"
$ws.Range("K6").Value  = " elena.dadyukova@yandex.ru"
$ws.Range("K7").Value  = "nataliya-klubnichka@mail.ru"
$ws.Range("K8").Value  = "Aleksandr541@inbox.ru"
$ws.Range("K9").Value  = "gulya2182@gmail.com"
$ws.Range("K12").Value = "-"
$ws.Range("K13").Value = "-"
$ws.Range("K14").Value = "karlinskay@mail.ru"
$ws.Range("K16").Value = "-"

# Родители cell that was previously empty for Краснов Кирилл Александрович
$ws.Range("D13").Value = "-"

# E-mail cells that the author turned into clickable mailto: hyperlinks
$ws.Hyperlinks.Add($ws.Range("K5"), "mailto:tatiazam@mail.ru", "", "", "tatiazam@mail.ru")
$ws.Hyperlinks.Add($ws.Range("K10"), "mailto:asemenova077@gmail.com", "", "", "asemenova077@gmail.com")
$ws.Hyperlinks.Add($ws.Range("K11"), "mailto:ingazarubina1983@gmail.com", "", "", "ingazarubina1983@gmail.com")
$ws.Hyperlinks.Add($ws.Range("K17"), "mailto:mcheblukov@bk.ru", "", "", "mcheblukov@bk.ru")

# Cosmetic: column E was manually narrowed to a width of 12 by the author
$ws.Columns.Item(5).ColumnWidth = 11.14

# Cosmetic: last selected cell before saving
$null = $ws.Range("Q11").Select()
